# Update the "想去人数" (want-to-go count) values in column F for the
# rows that changed between the two scraped snapshots.
# Same five events are listed on both the "展览" (Exhibition) sheet and the
# "全部类型" (All types) aggregate sheet, but at different row numbers there.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"     = @{ "F2" = 4657; "F3" = 2530; "F7"  = 62; "F12" = 1739; "F14" = 3886 }
    "全部类型" = @{ "F2" = 4657; "F3" = 2530; "F8"  = 62; "F16" = 1739; "F18" = 3886 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
